# edit.ps1 - applies the JsonForField.pptx edit described by the commit:
#   - fix something in input field (widen the JSON card + extend JSON text
#     with "cssField" / "scriptField" entries)
#   - fixed css field in json
#   - fixed script field in json
#
# Also re-caches the "datetimeFigureOut" field shown on the slide master and
# every slide layout from 1/1/2020 -> 1/2/2020 (a side effect the original
# author's PowerPoint produced when it re-saved the deck).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date field text (slide master + all 11 layouts)
# ---------------------------------------------------------------------
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "1/1/2020") {
            $tr.Text = "1/2/2020"
        }
    }
}

$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    for ($si = 1; $si -le $cl.Shapes.Count; $si++) {
        $sh = $cl.Shapes.Item($si)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "1/1/2020") {
                $tr.Text = "1/2/2020"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 2 - "Rounded Rectangle 4" JSON card
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$card = $s2.Shapes.Item(2)

# Widen the card to make room for the two extra JSON lines.
$card.Width = 366.909134

$tr = $card.TextFrame.TextRange

# Rewrite the body text - this recreates the paragraph layout (splitting on
# the PowerPoint paragraph-mark `vbCr`), then each paragraph below is
# re-split into the same multi-run shape the original file used.
$tr.Text = "{`r`t`"field_type`":`"input`",`r`t`"id`":`"ID`",`r`t`"type`":`"password`"`r`t`"place_holder`":`"AAA`",`r`t`"form_control_name`":`"abc`",`r`t`"cssField`":`"XXXXX`",`r`t`"scriptField`":`"XXXXX`"`r}"

# --- Paragraph 2: `"field_type":"input",  -> split into 3 runs ---
$c = $tr.Characters(3, 2);  $c.Text = $c.Text
$c = $tr.Characters(5, 18); $c.Text = $c.Text

# --- Paragraph 3: `"id":"ID",  -> split into 3 runs ---
$c = $tr.Characters(26, 2); $c.Text = $c.Text
$c = $tr.Characters(28, 7); $c.Text = $c.Text

# --- Paragraph 4: `"type":"password"  -> split into 3 runs ---
$c = $tr.Characters(38, 2);  $c.Text = $c.Text
$c = $tr.Characters(40, 15); $c.Text = $c.Text

# --- Paragraph 5: `"place_holder":"AAA",  -> split into 3 runs ---
$c = $tr.Characters(57, 2);  $c.Text = $c.Text
$c = $tr.Characters(59, 18); $c.Text = $c.Text

# --- Paragraph 6: `"form_control_name":"abc",  -> split into 3 runs ---
$c = $tr.Characters(80, 22);  $c.Text = $c.Text
$c = $tr.Characters(102, 3);  $c.Text = $c.Text

# --- Paragraph 7 (new): `"cssField":"XXXXX",  -> split into 3 runs ---
$c = $tr.Characters(108, 2); $c.Text = $c.Text
$c = $tr.Characters(110, 8); $c.Text = $c.Text

# --- Paragraph 8 (new): `"scriptField":"XXXXX"  -> split into 4 runs ---
$c = $tr.Characters(129, 1);  $c.Text = $c.Text
$c = $tr.Characters(130, 1);  $c.Text = $c.Text
$c = $tr.Characters(131, 11); $c.Text = $c.Text
